# Adding an extra layer of vegetarian/vegan substitutions to each red-meat
# sub-type sheet (pork, goat, beef, lamb, bison, venison, rabbit): a new
# row 5 with "vegetarian substitute" in column A and the specific
# substitute food in column B.

$wb = $excel.ActiveWorkbook

# Touch the "meat" sheet's selection first so it doesn't end up as the
# last-activated (tabSelected) sheet.
$ws = $wb.Worksheets.Item("meat")
$ws.Range("A11").Select()

# pork -> plant-based pork
$ws = $wb.Worksheets.Item("pork")
$ws.Range("A5").Value = "vegetarian substitute"
$ws.Range("B5").Value = "plant-based pork"
$ws.Range("A5:B5").Select()

# goat -> tofu
$ws = $wb.Worksheets.Item("goat")
$ws.Range("A5").Value = "vegetarian substitute"
$ws.Range("B5").Value = "tofu"
$ws.Range("B6").Select()

# beef -> plant-based ground beef
$ws = $wb.Worksheets.Item("beef")
$ws.Range("A5").Value = "vegetarian substitute"
$ws.Range("B5").Value = "plant-based ground beef"
$ws.Range("A5:B5").Select()

# lamb -> plant-based ground beef
$ws = $wb.Worksheets.Item("lamb")
$ws.Range("A5").Value = "vegetarian substitute"
$ws.Range("B5").Value = "plant-based ground beef"
$ws.Range("A5:B5").Select()

# bison -> plant-based ground beef
$ws = $wb.Worksheets.Item("bison")
$ws.Range("A5").Value = "vegetarian substitute"
$ws.Range("B5").Value = "plant-based ground beef"
$ws.Range("A5:B5").Select()

# venison -> tofu
$ws = $wb.Worksheets.Item("venison")
$ws.Range("A5").Value = "vegetarian substitute"
$ws.Range("B5").Value = "tofu"
$ws.Range("B6").Select()

# rabbit -> tofu (this sheet ends up active/selected, matching the target)
$ws = $wb.Worksheets.Item("rabbit")
$ws.Range("A5").Value = "vegetarian substitute"
$ws.Range("B5").Value = "tofu"
$ws.Range("A6").Select()
